# Delete the "NUMBER OF PROGRAMME RUNS" column (column H), shifting
# "SCRAP" (column I) left into column H.
$ws = $excel.ActiveWorkbook.ActiveSheet
$ws.Columns.Item(8).Delete()

# Delete row 7 (the "test.LST" entry).
$ws.Rows.Item(7).Delete()
